$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price (D) and Volume (E) columns retain their text formatting
# (values like "6.88" or "1.00" would otherwise be auto-converted to numbers).
$rng = $ws.Range("D2:E51")
$rng.NumberFormat = "@"

$ws.Range("D2").Value = '69.273.30'
$ws.Range("E2").Value = '  +0.01%  '
$ws.Range("D3").Value = '3.669.75'
$ws.Range("E3").Value = '  -0.29%  '
$ws.Range("D5").Value = '675.14'
$ws.Range("E5").Value = '  -1.19%  '
$ws.Range("D6").Value = '158.07'
$ws.Range("E6").Value = '  -2.44%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -1.26%  '
$ws.Range("E9").Value = '  -1.87%  '
$ws.Range("D10").Value = '6.88'
$ws.Range("E10").Value = '  -6.09%  '
$ws.Range("E11").Value = '  -2.51%  '
$ws.Range("E12").Value = '  -3.83%  '
$ws.Range("D13").Value = '4.289.04'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = '32.31'
$ws.Range("E14").Value = '  -3.89%  '
$ws.Range("D15").Value = '3.665.89'
$ws.Range("E15").Value = '  -0.36%  '
$ws.Range("D16").Value = '69.183.50'
$ws.Range("E16").Value = '  -0.26%  '
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("D18").Value = '15.99'
$ws.Range("E18").Value = '  -2.06%  '
$ws.Range("E19").Value = '  -2.98%  '
$ws.Range("D20").Value = '466.45'
$ws.Range("E20").Value = '  -2.90%  '
$ws.Range("D21").Value = '9.95'
$ws.Range("E21").Value = '  +0.37%  '
$ws.Range("D22").Value = '0.647'
$ws.Range("E22").Value = '  -2.77%  '
$ws.Range("D23").Value = '79.73'
$ws.Range("E23").Value = '  -0.57%  '
$ws.Range("D24").Value = '3.817.39'
$ws.Range("E24").Value = '  -0.26%  '
$ws.Range("E25").Value = '  -0.05%  '
$ws.Range("E26").Value = '  -6.72%  '
$ws.Range("E27").Value = '  -5.34%  '
$ws.Range("D28").Value = '9.03'
$ws.Range("E28").Value = '  -4.90%  '
$ws.Range("E29").Value = '  -1.38%  '
$ws.Range("E30").Value = '  -4.68%  '
$ws.Range("E31").Value = '  -3.43%  '
$ws.Range("D32").Value = '1.00'
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("E33").Value = '  -0.81%  '
$ws.Range("E34").Value = '  -4.88%  '
$ws.Range("D35").Value = '3.664.45'
$ws.Range("E35").Value = '  +0.39%  '
$ws.Range("D36").Value = '0.160'
$ws.Range("E36").Value = '  -4.71%  '
$ws.Range("D37").Value = '8.15'
$ws.Range("E37").Value = '  -3.91%  '
$ws.Range("D38").Value = '6.22'
$ws.Range("E38").Value = '  -1.85%  '
$ws.Range("E39").Value = '  +0.01%  '
$ws.Range("E40").Value = '  -0.12%  '
$ws.Range("E41").Value = '  -3.62%  '
$ws.Range("D42").Value = '0.0898'
$ws.Range("E42").Value = '  -3.88%  '
$ws.Range("D43").Value = '172.96'
$ws.Range("E43").Value = '  +6.53%  '
$ws.Range("D44").Value = '0.939'
$ws.Range("E44").Value = '  -1.38%  '
$ws.Range("E45").Value = '  -1.34%  '
$ws.Range("D46").Value = '28.23'
$ws.Range("E46").Value = '  -5.44%  '
$ws.Range("D47").Value = '0.000276'
$ws.Range("E47").Value = '  -3.59%  '
$ws.Range("D48").Value = '2.67'
$ws.Range("E48").Value = '  -5.18%  '
$ws.Range("D49").Value = '1.28'
$ws.Range("E49").Value = '  -4.48%  '
$ws.Range("E50").Value = '  -3.93%  '
$ws.Range("D51").Value = '7.76'
$ws.Range("E51").Value = '  -3.07%  '

# Restore the original (default) cell style now that the values are set as text.
$rng.Style = "Normal"
